# Weekly update: insert a new daily price record as the first data row
# (row 22) for "Agrícola del Norte S.A. de Arica - Durazno", pushing the
# existing data rows (22-53) down by one (to 23-54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 22. Excel shifts every
# row below it down by one and extends the used range accordingly.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$row = 22
$ws.Cells.Item($row, 1).Value2  = 1
$ws.Cells.Item($row, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value2  = 44874
$ws.Cells.Item($row, 5).Value2  = 15
$ws.Cells.Item($row, 6).Value2  = "Fruta"
$ws.Cells.Item($row, 7).Value2  = 100103
$ws.Cells.Item($row, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value2  = 100103004
$ws.Cells.Item($row, 10).Value2 = "Durazno"
$ws.Cells.Item($row, 11).Value2 = "Florida King"
$ws.Cells.Item($row, 12).Value2 = "Segunda"
$ws.Cells.Item($row, 13).Value2 = 200
$ws.Cells.Item($row, 14).Value2 = 24000
$ws.Cells.Item($row, 15).Value2 = 25000
$ws.Cells.Item($row, 16).Value2 = 24500
$ws.Cells.Item($row, 17).Value2 = "$/bandeja 10 kilos granel"
$ws.Cells.Item($row, 18).Value2 = "Región de Coquimbo"
$ws.Cells.Item($row, 19).Value2 = 2450
$ws.Cells.Item($row, 20).Value2 = 10
